$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new blank rows before row 15. This pushes the existing rows
# 15-17 (the "Cultivar IV Región" / bandeja data) down to rows 20-22,
# preserving their content unchanged, and leaves rows 15-19 blank
# (inheriting formatting from the row below, including the date style
# on column D) ready to be populated with the new weekly records.
$ws.Range("A15:A19").EntireRow.Insert()

# Row 15: Sin especificar / Extra
$ws.Cells.Item(15, 1).Value = 3
$ws.Cells.Item(15, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(15, 3).Value = "Coquimbo"
$ws.Cells.Item(15, 4).Value = 45131
$ws.Cells.Item(15, 5).Value = 5
$ws.Cells.Item(15, 6).Value = 100112043
$ws.Cells.Item(15, 7).Value = "Pepino dulce"
$ws.Cells.Item(15, 8).Value = "Sin especificar"
$ws.Cells.Item(15, 9).Value = "Extra"
$ws.Cells.Item(15, 10).Value = 56
$ws.Cells.Item(15, 11).Value = 20000
$ws.Cells.Item(15, 12).Value = 20000
$ws.Cells.Item(15, 13).Value = 20000
$ws.Cells.Item(15, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(15, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(15, 16).Value = 1333
$ws.Cells.Item(15, 17).Value = 15
$ws.Cells.Item(15, 18).Value = "Hortaliza"

# Row 16: Sin especificar / Primera
$ws.Cells.Item(16, 1).Value = 3
$ws.Cells.Item(16, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(16, 3).Value = "Coquimbo"
$ws.Cells.Item(16, 4).Value = 45131
$ws.Cells.Item(16, 5).Value = 5
$ws.Cells.Item(16, 6).Value = 100112043
$ws.Cells.Item(16, 7).Value = "Pepino dulce"
$ws.Cells.Item(16, 8).Value = "Sin especificar"
$ws.Cells.Item(16, 9).Value = "Primera"
$ws.Cells.Item(16, 10).Value = 75
$ws.Cells.Item(16, 11).Value = 18000
$ws.Cells.Item(16, 12).Value = 18000
$ws.Cells.Item(16, 13).Value = 18000
$ws.Cells.Item(16, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(16, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(16, 16).Value = 1200
$ws.Cells.Item(16, 17).Value = 15
$ws.Cells.Item(16, 18).Value = "Hortaliza"

# Row 17: Sin especificar / Primera / Provincia de Petorca
$ws.Cells.Item(17, 1).Value = 3
$ws.Cells.Item(17, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(17, 3).Value = "Coquimbo"
$ws.Cells.Item(17, 4).Value = 45131
$ws.Cells.Item(17, 5).Value = 5
$ws.Cells.Item(17, 6).Value = 100112043
$ws.Cells.Item(17, 7).Value = "Pepino dulce"
$ws.Cells.Item(17, 8).Value = "Sin especificar"
$ws.Cells.Item(17, 9).Value = "Primera"
$ws.Cells.Item(17, 10).Value = 56
$ws.Cells.Item(17, 11).Value = 17000
$ws.Cells.Item(17, 12).Value = 17000
$ws.Cells.Item(17, 13).Value = 17000
$ws.Cells.Item(17, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(17, 15).Value = "Provincia de Petorca"
$ws.Cells.Item(17, 16).Value = 1133
$ws.Cells.Item(17, 17).Value = 15
$ws.Cells.Item(17, 18).Value = "Hortaliza"

# Row 18: Sin especificar / Segunda / Provincia de Limarí
$ws.Cells.Item(18, 1).Value = 3
$ws.Cells.Item(18, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(18, 3).Value = "Coquimbo"
$ws.Cells.Item(18, 4).Value = 45131
$ws.Cells.Item(18, 5).Value = 5
$ws.Cells.Item(18, 6).Value = 100112043
$ws.Cells.Item(18, 7).Value = "Pepino dulce"
$ws.Cells.Item(18, 8).Value = "Sin especificar"
$ws.Cells.Item(18, 9).Value = "Segunda"
$ws.Cells.Item(18, 10).Value = 56
$ws.Cells.Item(18, 11).Value = 13000
$ws.Cells.Item(18, 12).Value = 13000
$ws.Cells.Item(18, 13).Value = 13000
$ws.Cells.Item(18, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(18, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(18, 16).Value = 867
$ws.Cells.Item(18, 17).Value = 15
$ws.Cells.Item(18, 18).Value = "Hortaliza"

# Row 19: Sin especificar / Segunda / Provincia de Petorca
$ws.Cells.Item(19, 1).Value = 3
$ws.Cells.Item(19, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(19, 3).Value = "Coquimbo"
$ws.Cells.Item(19, 4).Value = 45131
$ws.Cells.Item(19, 5).Value = 5
$ws.Cells.Item(19, 6).Value = 100112043
$ws.Cells.Item(19, 7).Value = "Pepino dulce"
$ws.Cells.Item(19, 8).Value = "Sin especificar"
$ws.Cells.Item(19, 9).Value = "Segunda"
$ws.Cells.Item(19, 10).Value = 50
$ws.Cells.Item(19, 11).Value = 12000
$ws.Cells.Item(19, 12).Value = 12000
$ws.Cells.Item(19, 13).Value = 12000
$ws.Cells.Item(19, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(19, 15).Value = "Provincia de Petorca"
$ws.Cells.Item(19, 16).Value = 800
$ws.Cells.Item(19, 17).Value = 15
$ws.Cells.Item(19, 18).Value = "Hortaliza"

# Keep row D-column date number format (style index 2, the same one used
# throughout column D) consistent on the new rows
$ws.Range("D15:D19").NumberFormat = $ws.Range("D20").NumberFormat
